# Fix Training Data Issue (#48)
# The "Date" column (BF) held the source filename-derived label
# "5-3-2011-12" for every data row. That label actually corresponds to
# game date 2012-05-03 (NBA stats were shown one day off), so replace the
# text in BF2:BF31 with the corrected date string "2012-05-03".
#
# We must keep these as literal text cells (not auto-converted to Excel
# date serials) and must not leave any residual number-format/style
# applied to the cells, so we temporarily mark the range as Text before
# writing the values, then clear the formatting back to the original
# (default/general) state once the literal text is safely stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")

# Force text storage so "2012-05-03" isn't reinterpreted as a date serial.
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "2012-05-03"
}

# Remove the temporary Text formatting so the cells end up unstyled, just
# like they were before the edit.
$dateRange.ClearFormats()
